# "dodanie zapisu do plikow trajektorii" - Dziennik praktyk edits
$d = $word.ActiveDocument

# 1) Total hours worked: 184 -> 160
$d.Content.Find.Execute("184", $false, $true, $false, $false, $false, $true, 1, $false, "160", 2) | Out-Null

# 2) Day 1 (01.07.2020) log entry: note the tools used, joining the first two
#    sentences (the <w:br/> between them is removed) and adding the
#    "(Trello, Discord, Visual Studio, GitHub)." sentence.
$vtab = [char]11
$find2 = "przygotowanie wykorzystywanych narzędzi." + $vtab + "Zebranie"
$replace2 = "przygotowanie wykorzystywanych narzędzi. (Trello, Discord, Visual Studio, GitHub). Zebranie"
$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $replace2, 2) | Out-Null

# 3) Day 6 (06.07.2020) log entry: elaborate on the libraries used and on the
#    purpose of the Harris/SIFT algorithms.
$find3 = "Wyszukanie bibliotek do implementacji poszczególnych algorytmów."
$replace3 = "Wyszukanie bibliotek do implementacji poszczególnych algorytmów (biblioteki: Eigen, OpenCV)."
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $replace3, 2) | Out-Null

$find4 = "Implementacja algorytmów „Harris Corner Detector” oraz „SIFT”."
$replace4 = "Implementacja algorytmów „Harris Corner Detector” oraz „SIFT” w celu wykrycia punktów charakterystycznych w obrazie."
$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $replace4, 2) | Out-Null

# 4) Daily hours worked column: the first day is now 6h, every other
#    logged day becomes 7h (instead of 8h for every day). Edited cell by
#    cell (rather than a global Find/Replace) so the "8" inside dates such
#    as "8.07.2020" is left untouched.
$tbl = $d.Tables.Item(1)
$rowCount = $tbl.Rows.Count
for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $tbl.Cell($r, 3)
    if ($r -eq 2) {
        $cell.Range.Text = "6"
    } else {
        $cell.Range.Text = "7"
    }
}
